$wb = $excel.ActiveWorkbook

# Update the Overview sheet: file 89ec4140-7020-4012-9fe1-624c2b8a2ebb has now
# been handed back (in sync) for both locales instead of merely "Ready for handoff"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# Update the zh-cn detail sheet for the same file: status flips to handed back,
# and the handback is now timestamped
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-21 18:47:53"

# Update the de-de detail sheet for the same file: status flips to handed back,
# and the handback is now timestamped
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-21 18:47:59"
